# Update the rpc-reply message-id UUIDs embedded in the logged XML text
# stored in column J (rows 2-5) of the "logs" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("logs")

$replacements = @(
    @{ Cell = "J2"; Old = "e97bf60f-3e4b-4d13-9bb8-52a57e15824e"; New = "7dc29a4a-195a-4e3a-b5fd-9d3528afee08" },
    @{ Cell = "J3"; Old = "46ab6473-6fbe-44da-be2e-1a95719f1517"; New = "3c13dc21-94ad-4cf0-a2b9-4da9dc00efcd" },
    @{ Cell = "J4"; Old = "65ca91d0-06b5-4aa4-b2e8-8d8184255e65"; New = "bb320285-691f-4e73-a08b-5ab5727f8a58" },
    @{ Cell = "J5"; Old = "b11e67a3-bddb-42b5-bcc1-1e9d4ebb6ed2"; New = "bf787d04-6c1c-41b1-84cf-88218c3f3119" }
)

foreach ($r in $replacements) {
    $cell = $ws.Range($r.Cell)
    $text = $cell.Value2
    $text = $text.Replace($r.Old, $r.New)
    $cell.Value2 = $text
}
